# Turn the plain A1:A5 number list on Sheet1 into a bordered "roll / marks /
# pass or fail" student table living in C3:E12, with a "total" row at the
# bottom (C12:E12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old A1:A5 values so they don't linger outside the new table.
$ws.Range("A1:A5").ClearContents()

# Header row.
$ws.Range("C3").Value = "roll "
$ws.Range("D3").Value = "marks"
$ws.Range("E3").Value = "pass or fail"

# Footer row label ("total"); D12/E12 stay blank but are part of the table.
$ws.Range("C12").Value = "total"

# Give every cell of the table (header + 8 blank data rows + total row) a
# thin box border on all four sides.
$ws.Range("C3:E12").Borders.LineStyle = 1

# Match the author's final selection/active cell.
$ws.Range("C12").Select()

Write-Output "table written"
